# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect freshly generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 7928
$ws1.Range("F10").Value = 462
$ws1.Range("F17").Value = 5821
$ws1.Range("F20").Value = 1691
$ws1.Range("F22").Value = 375

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 7928
$ws4.Range("F10").Value = 462
$ws4.Range("F18").Value = 5821
$ws4.Range("F22").Value = 1691
$ws4.Range("F24").Value = 375
